$wb = $excel.ActiveWorkbook

# Cell value updates (sheet -> cell -> new value), derived from the
# Behemoth_Profits.xlsx OOXML diff (scheduled-runner market-price refresh).
$updates = @{
    "ALC" = @{
        "H17" = 1675.2307
        "J17" = 1675.2307
        "L17" = 5025.6921
        "N17" = -5361.6921
        "H28" = 900.5333000000001
        "I28" = 653.3333
        "K28" = 653.3333
        "M28" = -168.3333
        "H69" = 15099.6
        "I69" = 9000
        "J69" = 15777.333
        "K69" = 27000
        "L69" = 47331.999
        "M69" = -26126
        "N69" = -49079.999
        "H72" = 15099.6
        "I72" = 9000
        "J72" = 15777.333
        "K72" = 81000
        "L72" = 141995.997
        "M72" = -76632
        "N72" = -150731.997
        "H76" = 5586.706
        "I76" = 4284.857
        "J76" = 6498
        "K76" = 4284.857
        "L76" = 6498
        "M76" = -3969.857
        "N76" = -7128
        "H79" = 5586.706
        "I79" = 4284.857
        "J79" = 6498
        "K79" = 4284.857
        "L79" = 6498
        "M79" = -3192.857
        "N79" = -8682
        "I86" = 4495.5
        "J86" = 5964.143
        "K86" = 4495.5
        "L86" = 5964.143
        "M86" = -3372.5
        "N86" = -8210.143
        "I89" = 4495.5
        "J89" = 5964.143
        "K89" = 22477.5
        "L89" = 29820.715
        "M89" = -16861.5
        "N89" = -41052.715
        "H111" = 799.6667
        "I111" = 799.6667
        "J111" = 0
        "K111" = 2399.0001
        "L111" = 0
        "M111" = 667.9998999999998
        "H132" = 1531.3684
        "I132" = 805.45715
        "J132" = 10000.333
        "K132" = 2416.37145
        "L132" = 30000.999
        "M132" = 113.6285500000004
        "N132" = -35060.999
    }
    "ARM" = @{
        "H61" = 34096400
        "I61" = 27781864
        "K61" = 27781864
        "M61" = -27781652
        "H88" = 2445.842
        "I88" = 1703.5
        "J88" = 3270.6667
        "K88" = 1703.5
        "L88" = 3270.6667
        "M88" = -1297.5
        "N88" = -4082.6667
        "H91" = 2445.842
        "I91" = 1703.5
        "J91" = 3270.6667
        "K91" = 1703.5
        "L91" = 3270.6667
        "M91" = -299.5
        "N91" = -6078.6667
        "H102" = 5447.619
        "I102" = 5669.85
        "J102" = 1003
        "K102" = 5669.85
        "L102" = 1003
        "M102" = -4047.85
        "N102" = -4247
        "H132" = 4116.5347
        "J132" = 22166.5
        "L132" = 66499.5
        "N132" = -71559.5
        "H136" = 34096400
        "I136" = 27781864
        "K136" = 83345592
        "M136" = -83343042
    }
    "BSM" = @{
        "H10" = 5
        "J10" = 0
        "L10" = 0
        "H20" = 5674.727
        "J20" = 5900.3076
        "L20" = 5900.3076
        "N20" = -6394.3076
        "H86" = 3664.7058
        "I86" = 3567
        "J86" = 4397.5
        "K86" = 3567
        "L86" = 4397.5
        "M86" = -2444
        "N86" = -6643.5
        "H89" = 3664.7058
        "I89" = 3567
        "J89" = 4397.5
        "K89" = 17835
        "L89" = 21987.5
        "M89" = -12219
        "N89" = -33219.5
        "H105" = 1863.5778
        "I105" = 1190.3334
        "K105" = 1190.3334
    }
    "CRP" = @{
        "H119" = 50493.668
        "J119" = 50493.668
        "L119" = 50493.668
        "N119" = -60169.668
        "H122" = 1488.95
        "J122" = 1833.8
        "L122" = 5501.4
        "N122" = -10401.4
        "H141" = 356085.6
        "J141" = 384650.66
        "L141" = 384650.66
        "N141" = -395010.66
    }
    "CUL" = @{
        "H3" = 2331.375
        "I3" = 1330.6
        "J3" = 3999.3333
        "K3" = 3991.8
        "L3" = 11997.9999
        "M3" = -3879.8
        "N3" = -12221.9999
        "H131" = 15366.733
        "J131" = 16727.1
        "L131" = 50181.3
        "N131" = -60261.3
        "H134" = 11010.611
        "J134" = 12961.833
        "L134" = 38885.499
        "N134" = -49025.499
    }
    "GSM" = @{
        "H97" = 2109.2144
        "I97" = 2185.3
        "J97" = 1919
        "K97" = 2185.3
        "L97" = 1919
        "M97" = -1689.3
        "N97" = -2911
        "H126" = 3618.6956
        "I126" = 3020.625
        "K126" = 9061.875
        "M126" = -6591.875
        "H132" = 30305398
        "I132" = 37039280
        "K132" = 111117840
        "M132" = -111115310
    }
    "LTW" = @{
        "H68" = 4160.5
        "I68" = 3992.6
        "K68" = 3992.6
        "M68" = -3243.6
        "H71" = 4160.5
        "I71" = 3992.6
        "K71" = 19963
        "M71" = -16219
        "H122" = 5981.2354
        "I122" = 5299.4736
        "J122" = 6844.8
        "K122" = 15898.4208
        "L122" = 20534.4
        "M122" = -13448.4208
        "N122" = -25434.4
        "H127" = 49985
        "J127" = 49985
        "L127" = 49985
        "N127" = -59905
        "H132" = 1005712.75
        "I132" = 116092.78
        "K132" = 348278.34
        "M132" = -345748.34
    }
    "WVR" = @{
        "H96" = 2004.25
        "I96" = 1622.25
        "J96" = 2258.9167
        "K96" = 1622.25
        "L96" = 2258.9167
        "M96" = -249.25
        "N96" = -5004.9167
        "H124" = 165646.67
        "J124" = 165646.67
        "L124" = 165646.67
        "N124" = -175466.67
        "H129" = 99792.5
        "J129" = 99792.5
        "L129" = 99792.5
        "N129" = -109792.5
        "H132" = 405450.7
        "I132" = 4793.381
        "K132" = 14380.143
        "M132" = -11850.143
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

# Cells whose trailing (HQ-only-leve) profit column no longer applies
# after this refresh -- fully removed, not just blanked, per the diff.
$clears = @{
    "ALC" = @("N111")
    "BSM" = @("N10")
}

foreach ($sheetName in $clears.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $clears[$sheetName]) {
        $ws.Range($cellRef).ClearContents()
    }
}
